# The deck ships two DrawingML themes:
#   ppt/theme/theme1.xml -> "Office Theme" / "Office" colour scheme  (used by the Notes Master)
#   ppt/theme/theme2.xml -> "Integral"     / "Red Violet" colour scheme (used by the Slide Master,
#                                                                         i.e. the live presentation theme)
# The authored change swaps the two themes' content: the Slide Master's theme becomes the
# plain "Office" colour scheme, while the Notes Master's theme becomes "Red Violet"/"Integral".
#
# The live, rendering theme (the one every Slide/Master/ColorScheme/ThemeColorScheme object in
# the PowerPoint object model resolves to) is the Slide Master's theme (theme2.xml). Re-point it
# at the standard Office palette so the active design matches the target "Office Theme" colours.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function ToRGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# ThemeColorScheme index order: 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink.
$tcs.Item(1).RGB  = ToRGB 0x00 0x00 0x00   # dk1      000000
$tcs.Item(2).RGB  = ToRGB 0xFF 0xFF 0xFF   # lt1      FFFFFF
$tcs.Item(3).RGB  = ToRGB 0x44 0x54 0x6A   # dk2      44546A
$tcs.Item(4).RGB  = ToRGB 0xE7 0xE6 0xE6   # lt2      E7E6E6
$tcs.Item(5).RGB  = ToRGB 0x5B 0x9B 0xD5   # accent1  5B9BD5
$tcs.Item(6).RGB  = ToRGB 0xED 0x7D 0x31   # accent2  ED7D31
$tcs.Item(7).RGB  = ToRGB 0xA5 0xA5 0xA5   # accent3  A5A5A5
$tcs.Item(8).RGB  = ToRGB 0xFF 0xC0 0x00   # accent4  FFC000
$tcs.Item(9).RGB  = ToRGB 0x44 0x72 0xC4   # accent5  4472C4
$tcs.Item(10).RGB = ToRGB 0x70 0xAD 0x47   # accent6  70AD47
$tcs.Item(11).RGB = ToRGB 0x05 0x63 0xC1   # hlink    0563C1
$tcs.Item(12).RGB = ToRGB 0x95 0x4F 0x72   # folHlink 954F72
